$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# 1. Update cell text content (write D11 before C10 so new shared strings are
#    appended in the order the target workbook expects: "Minh Chau @ Nguyen"
#    then "Ha My # Nguyen").
$ws2.Range("D11").Value = "Minh Châu @ Nguyễn"
$ws2.Range("C10").Value = "Hà My # Nguyễn"
$ws2.Range("E11").Value = "1 error field"

# 2. Fix up cell formatting so it matches the simplified style set:
#    - D11 drops its wrap-text styling (same formatting as C11).
#    - E2 / E11 drop their (redundant) explicit fill flag, reusing the
#      plain wrap-text style already used by C10.
$ws2.Range("C11").Copy()
$ws2.Range("D11").PasteSpecial(-4122)

$ws2.Range("C10").Copy()
$ws2.Range("E2").PasteSpecial(-4122)
$ws2.Range("E11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# 3. The rows no longer hold wrapped, multi-line text, so let them return to
#    the sheet's default height instead of the old fixed 51pt.
$ws2.Rows.Item(10).AutoFit()
$ws2.Rows.Item(11).AutoFit()

# 4. Column D is split out from the combined C:D width and given its own
#    (wider) width; column C keeps the original shared width.
$ws2.Columns.Item(4).ColumnWidth = 34.86

# 5. Performance Tracker becomes the active sheet/tab, with a new selection.
$ws2.Activate()
$ws2.Range("D11").Select()
